$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 106, shifting existing rows
# 106-119 down to 107-120 (new weekly price record for Alcachofa /
# Macroferia Regional de Talca).
$ws.Rows("106:106").Insert()

$ws.Range("A106").Value2 = 5
$ws.Range("B106").Value2 = "Macroferia Regional de Talca"
$ws.Range("C106").Value2 = "Maule"
$ws.Range("D106").Value2 = 45131
$ws.Range("D106").NumberFormat = $ws.Range("D107").NumberFormat
$ws.Range("E106").Value2 = 7
$ws.Range("F106").Value2 = 100112013
$ws.Range("G106").Value2 = "Alcachofa"
$ws.Range("H106").Value2 = "Madrigal"
$ws.Range("I106").Value2 = "Primera"
$ws.Range("J106").Value2 = 300
$ws.Range("K106").Value2 = 15000
$ws.Range("L106").Value2 = 15000
$ws.Range("M106").Value2 = 15000
$ws.Range("N106").Value2 = "`$/caja 40 unidades"
$ws.Range("O106").Value2 = "Provincia del Elquí"
$ws.Range("P106").Value2 = 375
$ws.Range("Q106").Value2 = 40
$ws.Range("R106").Value2 = "Hortaliza"
